$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.180.61'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.911.70'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5071'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3927'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09322'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.142'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.95'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.402'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.92%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.919.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.11%  '
$ws.Range('B14').Value = 'Solana'
$ws.Range('C14').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.94'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.331'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06622'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('E20').Value = '  +1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.227'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.236.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.327'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.598'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.136.03'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.16'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '158.10'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '127.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.106'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1075'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.661'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.724'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06684'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02436'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('B38').Value = 'Algorand'
$ws.Range('C38').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2210'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.91%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.243'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.286'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +8.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6524'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.57'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.025'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6141'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.724'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.289'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.026'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '122.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.189'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.77%  '
